$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in column C
$ws.Range("C2").Value = 30
$ws.Range("C3").Value = 31
$ws.Range("C4").Value = 22
$ws.Range("C5").Value = 33
$ws.Range("C6").Value = 18
$ws.Range("C7").Value = 21
$ws.Range("C8").Value = 23
$ws.Range("C9").Value = 30
$ws.Range("C11").Value = 29
$ws.Range("C12").Value = 35
$ws.Range("C14").Value = 19
$ws.Range("C15").Value = 24
$ws.Range("C16").Value = 26
$ws.Range("C17").Value = 24
$ws.Range("C18").Value = 23

# Update text values in column B
$ws.Range("B10").Value = "<see>"
$ws.Range("B13").Value = "<git>"
$ws.Range("B14").Value = "<up>"
$ws.Range("B18").Value = "<hen>"
